$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "bus" info to cell M4 (new shared string entry)
$ws.Range("M4").Value = "🚌مساء:1`n🚌صباحا: 7"

# Match style used by neighboring cells in row 4 (wrap text, no special border/fill)
$ws.Range("M4").WrapText = $true

# Update the active selection to P7
$ws.Range("P7").Select()
